$p = $ppt.ActivePresentation

# 1) Slide 6 table: switch the table's style to the built-in style
#    {BDA5F6D4-FEC8-42EA-AF9A-DAD1E15E9209} (was the custom
#    {68116B39-4B1D-427B-895E-26A34C154B87} style).
$tableSlide = $p.Slides.Item(6)
$tableShape = $tableSlide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{BDA5F6D4-FEC8-42EA-AF9A-DAD1E15E9209}")

# 2) Presentation design: switch the theme colours from "Integral" to the
#    default "Office Theme" palette.
$master = $p.SlideMaster
$colors = $master.Theme.ThemeColorScheme

$colors.Item(1).RGB  = 0x000000   # dk1       - 000000
$colors.Item(2).RGB  = 0xFFFFFF   # lt1       - FFFFFF
$colors.Item(3).RGB  = 0x6A5444   # dk2       - 44546A
$colors.Item(4).RGB  = 0xE6E6E7   # lt2       - E7E6E6
$colors.Item(5).RGB  = 0xD59B5B   # accent1   - 5B9BD5
$colors.Item(6).RGB  = 0x317DED   # accent2   - ED7D31
$colors.Item(7).RGB  = 0xA5A5A5   # accent3   - A5A5A5
$colors.Item(8).RGB  = 0x00C0FF   # accent4   - FFC000
$colors.Item(9).RGB  = 0xC47244   # accent5   - 4472C4
$colors.Item(10).RGB = 0x47AD70   # accent6   - 70AD47
$colors.Item(11).RGB = 0xC16305   # hlink     - 0563C1
$colors.Item(12).RGB = 0x724F95   # folHlink  - 954F72
